$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet / update header title for the new "through" date
$ws.Name = "Through 2022-05-16"
$ws.Range("B1").Value = "May 2022 (through May 16)"

# Row 3
$ws.Range("Q3").Value = 4
$ws.Range("V3").Value = 1

# Row 4
$ws.Range("G4").Value = 4

# Row 5
$ws.Range("L5").Value = 3

# Row 7
$ws.Range("Q7").Value = 3
$ws.Range("AA7").Value = 2

# Row 10
$ws.Range("AK10").Value = 1

# Row 14
$ws.Range("Q14").Value = 1

# Row 20
$ws.Range("AA20").Value = 1

# Row 21
$ws.Range("L21").Value = 1

# Row 24
$ws.Range("B24").Value = 2

# Row 28
$ws.Range("G28").Value = 1

# Row 32
$ws.Range("G32").Value = 1
$ws.Range("AF32").Value = 1

# Row 35
$ws.Range("Q35").Value = 2

# Row 40
$ws.Range("AA40").Value = 1

# Row 61
$ws.Range("Q61").Value = 1

# Row 63
$ws.Range("G63").Value = 4

# Row 70
$ws.Range("AA70").Value = 1

# Row 77
$ws.Range("V77").Value = 2

# Row 84
$ws.Range("B84").Value = 1
